# Updates a set of odds values on Sheet1 to match the latest Betfair Back/Lay
# snapshot for 2026-01-21 (see commit message: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value = 5.8
$ws.Range("Q2").Value = 1.89

$ws.Range("S4").Value = 1.37

$ws.Range("P5").Value = 2.56

$ws.Range("Q6").Value = 1.54
$ws.Range("T6").Value = 1.52
$ws.Range("X6").Value = 27
$ws.Range("AA6").Value = 50
$ws.Range("AB6").Value = 20
$ws.Range("AK6").Value = 980

$ws.Range("F8").Value = 3.75

$ws.Range("G9").Value = 1.98
$ws.Range("H9").Value = 3.7

$ws.Range("H10").Value = 2.96
$ws.Range("P10").Value = 2.06

$ws.Range("AO12").Value = 80

$ws.Range("G13").Value = 9.6
$ws.Range("AA13").Value = 13

$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 2.04
$ws.Range("X14").Value = 21

$ws.Range("K15").Value = 3.85
$ws.Range("Q15").Value = 2
$ws.Range("T15").Value = 1.84
$ws.Range("Y15").Value = 16.5
$ws.Range("AB15").Value = 8.800000000000001
$ws.Range("AM15").Value = 130

$ws.Range("F16").Value = 1.48
$ws.Range("N16").Value = 7.6
$ws.Range("P16").Value = 3.15
$ws.Range("Q16").Value = 1.42
$ws.Range("R16").Value = 1.89
$ws.Range("S16").Value = 2.04
$ws.Range("U16").Value = 2.5
$ws.Range("Y16").Value = 38
$ws.Range("AK16").Value = 14
$ws.Range("AN16").Value = 4.5

$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 9.4
$ws.Range("R17").Value = 1.97
$ws.Range("S17").Value = 1.94
$ws.Range("AH17").Value = 42

$ws.Range("H18").Value = 21
$ws.Range("K18").Value = 13.5
